$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: end time (C28) moves from 20:00 to 19:30 -> 19.5/24 day-fraction
$ws.Range("C28").Value = 19.5 / 24

# New note in E28 describing the day's work, matching the style used by
# other wrapped note cells (center/middle aligned + wrap text)
$ws.Range("E28").HorizontalAlignment = -4108
$ws.Range("E28").VerticalAlignment = -4108
$ws.Range("E28").WrapText = $true
$ws.Range("E28").Value = "Unity UI redesign + XML + scripts + 4 Free Boost + Initial Feats"

# Row grows to fit the two-line wrapped note
$ws.Rows.Item(28).RowHeight = 30

# Update the active selection left by the author
$ws.Range("E28").Select() | Out-Null
